{"js": "// Office.js (Word JavaScript API) script\n// Updates the worksheet date and all \"two-digit \u00f7 one-digit\" problems\n// to the new values per the target revision.\n\nconst replacements = [\n  [\"2024-11-10 Sunday\", \"2024-11-11 Monday\"],\n  [\"39\u00f73=\", \"87\u00f75=\"],\n  [\"48\u00f75=\", \"61\u00f72=\"],\n  [\"41\u00f75=\", \"22\u00f76=\"],\n  [\"30\u00f76=\", \"54\u00f76=\"],\n  [\"80\u00f76=\", \"69\u00f75=\"],\n  [\"42\u00f74=\", \"15\u00f77=\"],\n  [\"55\u00f76=\", \"76\u00f75=\"],\n  [\"20\u00f77=\", \"65\u00f79=\"],\n  [\"47\u00f72=\", \"68\u00f74=\"],\n  [\"82\u00f78=\", \"41\u00f76=\"],\n  [\"90\u00f79=\", \"38\u00f77=\"],\n  [\"42\u00f72=\", \"96\u00f78=\"],\n  [\"53\u00f78=\", \"48\u00f75=\"],\n  [\"55\u00f74=\", \"13\u00f77=\"],\n  [\"22\u00f77=\", \"16\u00f76=\"],\n  [\"83\u00f76=\", \"44\u00f76=\"],\n  [\"46\u00f77=\", \"37\u00f77=\"],\n  [\"75\u00f76=\", \"16\u00f75=\"],\n  [\"42\u00f73=\", \"37\u00f76=\"],\n  [\"24\u00f78=\", \"10\u00f72=\"],\n  [\"42\u00f75=\", \"16\u00f78=\"],\n  [\"44\u00f79=\", \"33\u00f73=\"],\n  [\"86\u00f72=\", \"65\u00f75=\"],\n  [\"33\u00f75=\", \"48\u00f78=\"],\n  [\"26\u00f77=\", \"23\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Updates the worksheet date and all \"two-digit \u00f7 one-digit\" problems\n# to the new values per the target revision.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-11-10 Sunday\", \"2024-11-11 Monday\"),\n  @(\"39\u00f73=\", \"87\u00f75=\"),\n  @(\"48\u00f75=\", \"61\u00f72=\"),\n  @(\"41\u00f75=\", \"22\u00f76=\"),\n  @(\"30\u00f76=\", \"54\u00f76=\"),\n  @(\"80\u00f76=\", \"69\u00f75=\"),\n  @(\"42\u00f74=\", \"15\u00f77=\"),\n  @(\"55\u00f76=\", \"76\u00f75=\"),\n  @(\"20\u00f77=\", \"65\u00f79=\"),\n  @(\"47\u00f72=\", \"68\u00f74=\"),\n  @(\"82\u00f78=\", \"41\u00f76=\"),\n  @(\"90\u00f79=\", \"38\u00f77=\"),\n  @(\"42\u00f72=\", \"96\u00f78=\"),\n  @(\"53\u00f78=\", \"48\u00f75=\"),\n  @(\"55\u00f74=\", \"13\u00f77=\"),\n  @(\"22\u00f77=\", \"16\u00f76=\"),\n  @(\"83\u00f76=\", \"44\u00f76=\"),\n  @(\"46\u00f77=\", \"37\u00f77=\"),\n  @(\"75\u00f76=\", \"16\u00f75=\"),\n  @(\"42\u00f73=\", \"37\u00f76=\"),\n  @(\"24\u00f78=\", \"10\u00f72=\"),\n  @(\"42\u00f75=\", \"16\u00f78=\"),\n  @(\"44\u00f79=\", \"33\u00f73=\"),\n  @(\"86\u00f72=\", \"65\u00f75=\"),\n  @(\"33\u00f75=\", \"48\u00f78=\"),\n  @(\"26\u00f77=\", \"23\u00f74=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $new\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n\n  $found = $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $find.Format, $find.Replacement.Text, 2)\n\n  if (-not $found) {\n    Write-Output \"WARNING: text not found: $old\"\n  }\n}\n"}
